$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 label changes from "PERSON" placeholder row index shift (still "PERSON")
$ws.Range("A1").Value = "PERSON"
$ws.Range("B1").Value = "A"
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"
$ws.Range("E1").Value = "D"

# Data rows 2-16: Person name stays the same per row, new pick columns B-E populated
$data = @(
    @("Andy Davis", "Brooks Koepka", "Alexander Noren", "Patrick Cantlay", "Lucas Bjerregaard"),
    @("Alex Evans", "Rickie Fowler", "Rafael Cabrera Bello", "Patrick Cantlay", "Adam Long"),
    @("Joel Malzer", "Justin Rose", "Ian Poulter", "Cameron Smith", "Devon Bling"),
    @("Bradlee Larson", "Rory McIlroy", "Kevin Kisner", "Shane Lowry", "Mike Weir"),
    @("Patrick Daugherty", "Rory McIlroy", "Xander Schauffele", "Li Haotong", "Lucas Bjerregaard"),
    @("John Ekola", "Rory McIlroy", "Xander Schauffele", "Patrick Cantlay", "Justin Harding"),
    @("Ryan Naughton", "Rickie Fowler", "Louis Oosthuizen", "Corey Conners", "Vijay Singh"),
    @("Kevin Evans", "Rory McIlroy", "Kevin Kisner", "Kiradech Aphibarnrat", "Sandy Lyle"),
    @("Dave Schunk", "Dustin Johnson", "Brandt Snedeker", "Cameron Smith", "Mike Weir"),
    @("Rick Laszewski", "Dustin Johnson", "Xander Schauffele", "Patrick Cantlay", "Lucas Bjerregaard"),
    @("John Griffin", "Tiger Woods", "Gary Woodland", "Patrick Cantlay", "Lucas Bjerregaard"),
    @("Matt Hanse", "Rory McIlroy", "Louis Oosthuizen", "Cameron Smith", "Lucas Bjerregaard"),
    @("Brandon Griffin", "Rory McIlroy", "Louis Oosthuizen", "Patrick Cantlay", "Justin Harding"),
    @("Reyanna Ekola", "Tiger Woods", "Xander Schauffele", "Eddie Pepperell", "Satoshi Kodaira"),
    @("Ryan Schunk", "Tiger Woods", "Matthew Fitzpatrick", "Kiradech Aphibarnrat", "Satoshi Kodaira")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Update selection to match saved view state
$ws.Range("D7").Select() | Out-Null
